$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert a new '2022-Q1' sheet positioned right before '总计'
# ---------------------------------------------------------------
$totalSheetAnchor = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetAnchor)
$newSheet.Name = "2022-Q1"

# NOTE: the sheet object returned by Add(Before:=...) shares its underlying
# anchor with the "before" argument that was passed in, so $totalSheetAnchor
# now actually refers to the freshly-inserted sheet, not to "总计" any more.
# Re-resolve "总计" by name to get a correct, un-aliased reference to it.
$totalSheet = $wb.Worksheets.Item("总计")

# Reference sheet to copy cell formatting from (same layout as other quarters)
$refSheet = $wb.Worksheets.Item("2021-Q4")

# Copy header-row format (bold + border + centered) from the reference sheet
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy column-A index-cell format (bold + border + centered) for the data rows
$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A18").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows (A: running index, B: fund code, C: fund name, D: fund scale,
#  E: total stock position, F: position ratio, G: held market value, H: position rank)
# row 2
$newSheet.Cells.Item(2,1).Value = 0
$c = $newSheet.Cells.Item(2,2)
$c.NumberFormat = "@"
$c.Value = "012588"
$newSheet.Cells.Item(2,3).Value = "南方港股通优势企业混合型证券投资基金A"
$c = $newSheet.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = "37.54"
$c = $newSheet.Cells.Item(2,5)
$c.NumberFormat = "@"
$c.Value = "71.00"
$c = $newSheet.Cells.Item(2,6)
$c.NumberFormat = "@"
$c.Value = "2.85"
$c = $newSheet.Cells.Item(2,7)
$c.NumberFormat = "@"
$c.Value = "1.0699"
$newSheet.Cells.Item(2,8).Value = 4

# row 3
$newSheet.Cells.Item(3,1).Value = 1
$c = $newSheet.Cells.Item(3,2)
$c.NumberFormat = "@"
$c.Value = "008513"
$newSheet.Cells.Item(3,3).Value = "南方宝丰混合A"
$c = $newSheet.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = "51.62"
$c = $newSheet.Cells.Item(3,5)
$c.NumberFormat = "@"
$c.Value = "21.13"
$c = $newSheet.Cells.Item(3,6)
$c.NumberFormat = "@"
$c.Value = "0.43"
$c = $newSheet.Cells.Item(3,7)
$c.NumberFormat = "@"
$c.Value = "0.2220"
$newSheet.Cells.Item(3,8).Value = 9

# row 4
$newSheet.Cells.Item(4,1).Value = 2
$c = $newSheet.Cells.Item(4,2)
$c.NumberFormat = "@"
$c.Value = "862001"
$newSheet.Cells.Item(4,3).Value = "光大阳光香港精选混合型集合资产管理计划（QDII）A 人民币"
$c = $newSheet.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = "3.91"
$c = $newSheet.Cells.Item(4,5)
$c.NumberFormat = "@"
$c.Value = "89.45"
$c = $newSheet.Cells.Item(4,6)
$c.NumberFormat = "@"
$c.Value = "5.44"
$c = $newSheet.Cells.Item(4,7)
$c.NumberFormat = "@"
$c.Value = "0.2127"
$newSheet.Cells.Item(4,8).Value = 5

# row 5
$newSheet.Cells.Item(5,1).Value = 3
$c = $newSheet.Cells.Item(5,2)
$c.NumberFormat = "@"
$c.Value = "862011"
$newSheet.Cells.Item(5,3).Value = "光大阳光香港精选混合型集合资产管理计划（QDII）A 美元"
$c = $newSheet.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "3.91"
$c = $newSheet.Cells.Item(5,5)
$c.NumberFormat = "@"
$c.Value = "89.45"
$c = $newSheet.Cells.Item(5,6)
$c.NumberFormat = "@"
$c.Value = "5.44"
$c = $newSheet.Cells.Item(5,7)
$c.NumberFormat = "@"
$c.Value = "0.2127"
$newSheet.Cells.Item(5,8).Value = 5

# row 6
$newSheet.Cells.Item(6,1).Value = 4
$c = $newSheet.Cells.Item(6,2)
$c.NumberFormat = "@"
$c.Value = "862012"
$newSheet.Cells.Item(6,3).Value = "光大阳光香港精选混合型集合资产管理计划（QDII）C 人民币"
$c = $newSheet.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "3.91"
$c = $newSheet.Cells.Item(6,5)
$c.NumberFormat = "@"
$c.Value = "89.45"
$c = $newSheet.Cells.Item(6,6)
$c.NumberFormat = "@"
$c.Value = "5.44"
$c = $newSheet.Cells.Item(6,7)
$c.NumberFormat = "@"
$c.Value = "0.2127"
$newSheet.Cells.Item(6,8).Value = 5

# row 7
$newSheet.Cells.Item(7,1).Value = 5
$c = $newSheet.Cells.Item(7,2)
$c.NumberFormat = "@"
$c.Value = "010010"
$newSheet.Cells.Item(7,3).Value = "国投瑞银港股通6个月定期开放股票"
$c = $newSheet.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = "8.09"
$c = $newSheet.Cells.Item(7,5)
$c.NumberFormat = "@"
$c.Value = "93.58"
$c = $newSheet.Cells.Item(7,6)
$c.NumberFormat = "@"
$c.Value = "2.61"
$c = $newSheet.Cells.Item(7,7)
$c.NumberFormat = "@"
$c.Value = "0.2111"
$newSheet.Cells.Item(7,8).Value = 10

# row 8
$newSheet.Cells.Item(8,1).Value = 6
$c = $newSheet.Cells.Item(8,2)
$c.NumberFormat = "@"
$c.Value = "013200"
$newSheet.Cells.Item(8,3).Value = "南方均衡优选一年持有期混合A"
$c = $newSheet.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "9.72"
$c = $newSheet.Cells.Item(8,5)
$c.NumberFormat = "@"
$c.Value = "33.40"
$c = $newSheet.Cells.Item(8,6)
$c.NumberFormat = "@"
$c.Value = "1.72"
$c = $newSheet.Cells.Item(8,7)
$c.NumberFormat = "@"
$c.Value = "0.1672"
$newSheet.Cells.Item(8,8).Value = 2

# row 9
$newSheet.Cells.Item(9,1).Value = 7
$c = $newSheet.Cells.Item(9,2)
$c.NumberFormat = "@"
$c.Value = "860007"
$newSheet.Cells.Item(9,3).Value = "光大阳光价值30个月持有期混合A"
$c = $newSheet.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "2.97"
$c = $newSheet.Cells.Item(9,5)
$c.NumberFormat = "@"
$c.Value = "90.90"
$c = $newSheet.Cells.Item(9,6)
$c.NumberFormat = "@"
$c.Value = "5.18"
$c = $newSheet.Cells.Item(9,7)
$c.NumberFormat = "@"
$c.Value = "0.1538"
$newSheet.Cells.Item(9,8).Value = 4

# row 10
$newSheet.Cells.Item(10,1).Value = 8
$c = $newSheet.Cells.Item(10,2)
$c.NumberFormat = "@"
$c.Value = "860027"
$newSheet.Cells.Item(10,3).Value = "光大阳光价值30个月持有期混合B"
$c = $newSheet.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "2.31"
$c = $newSheet.Cells.Item(10,5)
$c.NumberFormat = "@"
$c.Value = "90.90"
$c = $newSheet.Cells.Item(10,6)
$c.NumberFormat = "@"
$c.Value = "5.18"
$c = $newSheet.Cells.Item(10,7)
$c.NumberFormat = "@"
$c.Value = "0.1197"
$newSheet.Cells.Item(10,8).Value = 4

# row 11
$newSheet.Cells.Item(11,1).Value = 9
$c = $newSheet.Cells.Item(11,2)
$c.NumberFormat = "@"
$c.Value = "005741"
$newSheet.Cells.Item(11,3).Value = "南方君信灵活配置混合A"
$c = $newSheet.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "3.69"
$c = $newSheet.Cells.Item(11,5)
$c.NumberFormat = "@"
$c.Value = "76.52"
$c = $newSheet.Cells.Item(11,6)
$c.NumberFormat = "@"
$c.Value = "1.70"
$c = $newSheet.Cells.Item(11,7)
$c.NumberFormat = "@"
$c.Value = "0.0627"
$newSheet.Cells.Item(11,8).Value = 9

# row 12
$newSheet.Cells.Item(12,1).Value = 10
$c = $newSheet.Cells.Item(12,2)
$c.NumberFormat = "@"
$c.Value = "012589"
$newSheet.Cells.Item(12,3).Value = "南方港股通优势企业混合型证券投资基金C"
$c = $newSheet.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "2.05"
$c = $newSheet.Cells.Item(12,5)
$c.NumberFormat = "@"
$c.Value = "71.00"
$c = $newSheet.Cells.Item(12,6)
$c.NumberFormat = "@"
$c.Value = "2.85"
$c = $newSheet.Cells.Item(12,7)
$c.NumberFormat = "@"
$c.Value = "0.0584"
$newSheet.Cells.Item(12,8).Value = 4

# row 13
$newSheet.Cells.Item(13,1).Value = 11
$c = $newSheet.Cells.Item(13,2)
$c.NumberFormat = "@"
$c.Value = "008514"
$newSheet.Cells.Item(13,3).Value = "南方宝丰混合C"
$c = $newSheet.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "4.72"
$c = $newSheet.Cells.Item(13,5)
$c.NumberFormat = "@"
$c.Value = "21.13"
$c = $newSheet.Cells.Item(13,6)
$c.NumberFormat = "@"
$c.Value = "0.43"
$c = $newSheet.Cells.Item(13,7)
$c.NumberFormat = "@"
$c.Value = "0.0203"
$newSheet.Cells.Item(13,8).Value = 9

# row 14
$newSheet.Cells.Item(14,1).Value = 12
$c = $newSheet.Cells.Item(14,2)
$c.NumberFormat = "@"
$c.Value = "013201"
$newSheet.Cells.Item(14,3).Value = "南方均衡优选一年持有期混合C"
$c = $newSheet.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "0.84"
$c = $newSheet.Cells.Item(14,5)
$c.NumberFormat = "@"
$c.Value = "33.40"
$c = $newSheet.Cells.Item(14,6)
$c.NumberFormat = "@"
$c.Value = "1.72"
$c = $newSheet.Cells.Item(14,7)
$c.NumberFormat = "@"
$c.Value = "0.0144"
$newSheet.Cells.Item(14,8).Value = 2

# row 15
$newSheet.Cells.Item(15,1).Value = 13
$c = $newSheet.Cells.Item(15,2)
$c.NumberFormat = "@"
$c.Value = "501303"
$newSheet.Cells.Item(15,3).Value = "广发港股通恒生综合中型股指数(LOF)A"
$c = $newSheet.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "0.34"
$c = $newSheet.Cells.Item(15,5)
$c.NumberFormat = "@"
$c.Value = "92.39"
$c = $newSheet.Cells.Item(15,6)
$c.NumberFormat = "@"
$c.Value = "1.43"
$c = $newSheet.Cells.Item(15,7)
$c.NumberFormat = "@"
$c.Value = "0.0049"
$newSheet.Cells.Item(15,8).Value = 8

# row 16
$newSheet.Cells.Item(16,1).Value = 14
$c = $newSheet.Cells.Item(16,2)
$c.NumberFormat = "@"
$c.Value = "004996"
$newSheet.Cells.Item(16,3).Value = "广发港股通恒生综合中型股指数(LOF)C"
$c = $newSheet.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = "0.11"
$c = $newSheet.Cells.Item(16,5)
$c.NumberFormat = "@"
$c.Value = "92.39"
$c = $newSheet.Cells.Item(16,6)
$c.NumberFormat = "@"
$c.Value = "1.43"
$c = $newSheet.Cells.Item(16,7)
$c.NumberFormat = "@"
$c.Value = "0.0016"
$newSheet.Cells.Item(16,8).Value = 8

# row 17
$newSheet.Cells.Item(17,1).Value = 15
$c = $newSheet.Cells.Item(17,2)
$c.NumberFormat = "@"
$c.Value = "160922"
$newSheet.Cells.Item(17,3).Value = "大成恒生综合中小型股指数(QDII-LOF)A"
$c = $newSheet.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = "0.10"
$c = $newSheet.Cells.Item(17,5)
$c.NumberFormat = "@"
$c.Value = "92.44"
$c = $newSheet.Cells.Item(17,6)
$c.NumberFormat = "@"
$c.Value = "1.12"
$c = $newSheet.Cells.Item(17,7)
$c.NumberFormat = "@"
$c.Value = "0.0011"
$newSheet.Cells.Item(17,8).Value = 8

# row 18
$newSheet.Cells.Item(18,1).Value = 16
$c = $newSheet.Cells.Item(18,2)
$c.NumberFormat = "@"
$c.Value = "010150"
$newSheet.Cells.Item(18,3).Value = "南方君信灵活配置混合C"
$c = $newSheet.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = "0.02"
$c = $newSheet.Cells.Item(18,5)
$c.NumberFormat = "@"
$c.Value = "76.52"
$c = $newSheet.Cells.Item(18,6)
$c.NumberFormat = "@"
$c.Value = "1.70"
$c = $newSheet.Cells.Item(18,7)
$c.NumberFormat = "@"
$c.Value = "0.0003"
$newSheet.Cells.Item(18,8).Value = 9

# ---------------------------------------------------------------
# 2. Update the '总计' sheet: insert a row for 2022-Q1 at the top
#    of the data block, then renumber the running index column.
# ---------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 17
$totalSheet.Cells.Item(2,4).Value = 2.75

# Renumber the A-column running index (0,1,2,3,4,5) now that a row was added
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

